$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B2").Value = 209
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 66
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 23

# Remove the now-unused row 5 entirely (shrinks used range to A1:B4)
$ws.Range("A5:B5").Delete()
